$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Sat Oct 12 23:52:51 EDT 2024"
$ws.Range("B3").Value = "Sat Oct 12 23:53:03 EDT 2024"
$ws.Range("B4").Value = "Sat Oct 12 23:53:15 EDT 2024"
$ws.Range("B5").Value = "Sat Oct 12 23:53:27 EDT 2024"
$ws.Range("B6").Value = "Sat Oct 12 23:53:39 EDT 2024"
$ws.Range("B7").Value = "Sat Oct 12 23:53:51 EDT 2024"
$ws.Range("B8").Value = "Sat Oct 12 23:54:03 EDT 2024"
$ws.Range("B9").Value = "Sat Oct 12 23:54:16 EDT 2024"
$ws.Range("B10").Value = "Sat Oct 12 23:54:28 EDT 2024"
$ws.Range("B11").Value = "Sat Oct 12 23:54:40 EDT 2024"
$ws.Range("B12").Value = "Sat Oct 12 23:54:52 EDT 2024"
$ws.Range("B13").Value = "Sat Oct 12 23:55:04 EDT 2024"
$ws.Range("B14").Value = "Sat Oct 12 23:55:16 EDT 2024"
$ws.Range("B15").Value = "Sat Oct 12 23:55:29 EDT 2024"
$ws.Range("B16").Value = "Sat Oct 12 23:55:41 EDT 2024"
$ws.Range("B17").Value = "Sat Oct 12 23:55:54 EDT 2024"
$ws.Range("B18").Value = "Sat Oct 12 23:56:07 EDT 2024"
$ws.Range("B19").Value = "Sat Oct 12 23:56:20 EDT 2024"
$ws.Range("B20").Value = "Sat Oct 12 23:56:33 EDT 2024"
$ws.Range("B21").Value = "Sat Oct 12 23:56:47 EDT 2024"
$ws.Range("B22").Value = "Sat Oct 12 23:57:00 EDT 2024"
$ws.Range("B23").Value = "Sat Oct 12 23:57:13 EDT 2024"
$ws.Range("B24").Value = "Sat Oct 12 23:57:26 EDT 2024"
$ws.Range("B25").Value = "Sat Oct 12 23:57:39 EDT 2024"
$ws.Range("B26").Value = "Sat Oct 12 23:57:52 EDT 2024"
$ws.Range("B27").Value = "Sat Oct 12 23:58:05 EDT 2024"
$ws.Range("B28").Value = "Sat Oct 12 23:58:18 EDT 2024"
$ws.Range("B29").Value = "Sat Oct 12 23:58:31 EDT 2024"
$ws.Range("B30").Value = "Sat Oct 12 23:58:44 EDT 2024"
$ws.Range("B31").Value = "Sat Oct 12 23:58:57 EDT 2024"
$ws.Range("B32").Value = "Sat Oct 12 23:59:10 EDT 2024"
$ws.Range("B33").Value = "Sat Oct 12 23:59:23 EDT 2024"
$ws.Range("B34").Value = "Sat Oct 12 23:59:37 EDT 2024"
$ws.Range("B35").Value = "Sat Oct 12 23:59:49 EDT 2024"
$ws.Range("B36").Value = "Sun Oct 13 00:00:01 EDT 2024"
$ws.Range("B37").Value = "Sun Oct 13 00:00:13 EDT 2024"
$ws.Range("B38").Value = "Sun Oct 13 00:00:25 EDT 2024"
$ws.Range("B39").Value = "Sun Oct 13 00:00:37 EDT 2024"
$ws.Range("B40").Value = "Sun Oct 13 00:00:49 EDT 2024"
$ws.Range("B41").Value = "Sun Oct 13 00:01:00 EDT 2024"
$ws.Range("B42").Value = "Sun Oct 13 00:01:12 EDT 2024"
$ws.Range("B43").Value = "Sun Oct 13 00:01:26 EDT 2024"
$ws.Range("B44").Value = "Sun Oct 13 00:01:39 EDT 2024"
$ws.Range("B45").Value = "Sun Oct 13 00:01:52 EDT 2024"
$ws.Range("B46").Value = "Sun Oct 13 00:02:05 EDT 2024"
$ws.Range("B47").Value = "Sun Oct 13 00:02:18 EDT 2024"
$ws.Range("B48").Value = "Sun Oct 13 00:02:31 EDT 2024"
$ws.Range("B49").Value = "Sun Oct 13 00:02:44 EDT 2024"
$ws.Range("B50").Value = "Sun Oct 13 00:02:57 EDT 2024"
$ws.Range("B51").Value = "Sun Oct 13 00:03:10 EDT 2024"
$ws.Range("B52").Value = "Sun Oct 13 00:03:23 EDT 2024"
$ws.Range("B53").Value = "Sun Oct 13 00:03:36 EDT 2024"
$ws.Range("B54").Value = "Sun Oct 13 00:03:48 EDT 2024"
